$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Asia
$ws.Range("C3").Value = 217408107
$ws.Range("D3").Value = 22142
$ws.Range("E3").Value = 1546495
$ws.Range("F3").Value = 35
$ws.Range("G3").Value = 201093712
$ws.Range("H3").Value = 11705
$ws.Range("I3").Value = 14767900
$ws.Range("J3").Value = 15409

# Row 4 - Europe
$ws.Range("C4").Value = 249318460
$ws.Range("D4").Value = 503
$ws.Range("E4").Value = 2060520
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 245288281
$ws.Range("H4").Value = 4404
$ws.Range("I4").Value = 1969659

# Row 6 - Oceania
$ws.Range("C6").Value = 14323773
$ws.Range("D6").Value = 1
$ws.Range("I6").Value = 145455
